# 20200403{Start a new day_3}
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: new log entry "Half of it in the product cipher encryption"
$ws.Range("E5").Value = "Ok"
$ws.Range("F5").Value = "Half of it in the product cipher encryption"
$ws.Range("G5").Value = 0.34722222222222227
$ws.Range("G5").NumberFormat = "h:mm"

# Row 6: new log entry "A headache on my head and I can't forgive myself"
$ws.Range("E6").Value = "Ok"
$ws.Range("F6").Value = "A headache on my head and I can't forgive myself"
$ws.Range("G6").Value = 0.41666666666666669
$ws.Range("G6").NumberFormat = "h:mm"

# Row 7: new log entry "Alhamdu Li Allah => sleeping period"
$ws.Range("E7").Value = "Ok"
$ws.Range("F7").Value = "Alhamdu Li Allah => sleeping period"
$ws.Range("G7").Value = 0.16666666666666666
$ws.Range("G7").NumberFormat = "h:mm"

# Row 8: new log entry "Alhamdu Li Allah " (trailing space)
$ws.Range("E8").Value = "Ok"
$ws.Range("F8").Value = "Alhamdu Li Allah "
$ws.Range("G8").Value = 0.17361111111111113
$ws.Range("G8").NumberFormat = "h:mm"

# Row 10 (task list): rename "Programming Blazor" -> "Programming GP"
$ws.Range("B10").Value = "Programming GP"

# Re-create the A49:A68 / A14:A38 merged blocks so mergeCells ends up
# reordered the same way the original authoring session produced it:
# by touching (unmerge/merge) every OTHER merged block, those blocks get
# pushed to the end of the internal merge list, leaving A49:A68/A14:A38
# at the front, in their original relative order.
$rest = @(
    "A129:A133","A2:A13","A99:A103","A104:A108","A109:A113","A114:A118",
    "A119:A123","A124:A128","A69:A73","A74:A78","A79:A83","A84:A88",
    "A89:A93","A94:A98","A39:A43","A44:A48"
)
foreach ($r in $rest) {
    $ws.Range($r).UnMerge()
    $ws.Range($r).Merge()
}

# Update the active selection to F18
$ws.Range("F18").Select() | Out-Null
